# Apply the new block-order values for the "60_scenecat_block_order" sheet.
# Header row: swap the labels in A1/B1 (bedrooms_1 <-> living_rooms_1).
# Data rows 2-7: rewrite the one-hot "1" marker into its new column per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("D1").Value = "living_rooms_2"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "kitchens_2"

# Data rows (rows 2-7), one-hot encoded block order
$data = @(
    @(0,1,0,0,0,0),  # row 2
    @(0,0,0,0,0,1),  # row 3
    @(0,0,0,1,0,0),  # row 4
    @(0,0,0,0,1,0),  # row 5
    @(0,0,1,0,0,0),  # row 6
    @(1,0,0,0,0,0)   # row 7
)

$cols = @("A","B","C","D","E","F")

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $cellRef = "$($cols[$j])$row"
        $ws.Range($cellRef).Value = $values[$j]
    }
}
